$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected some selection scopes: clear out-of-range forecast values
# (columns G and H) for rows 3,5,7,9,11,13,15,17 where the forecast
# horizon extends beyond the available evaluation data.
$rows = @(3, 5, 7, 9, 11, 13, 15, 17)
foreach ($r in $rows) {
    $ws.Range("G${r}:H${r}").ClearContents()
}
